$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at the top of the "Repollo" block (row 1048) to make room
# for a new week's worth of price data (2021-09-09, serial 44448). This shifts
# the existing rows 1048:1067 down to 1054:1073, matching the new dimension
# A1:R1073.
$ws.Rows("1048:1053").Insert()

# New rows for date 44448 (2021-09-09), mirroring the existing row layout:
# Copenhague/Primera, Copenhague/Segunda, Crespo record/Primera,
# Crespo record/Segunda, Morada(o)/Primera, Morada(o)/Segunda.
$newRows = @(
  @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44448, 4, 100112006, 'Repollo', 'Copenhague', 'Primera', 2400, 850, 900, 875, '$/unidad', 'Provincia del Elquí', 875, 1, 'Hortaliza'),
  @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44448, 4, 100112006, 'Repollo', 'Copenhague', 'Segunda', 1320, 700, 800, 750, '$/unidad', 'Provincia del Elquí', 750, 1, 'Hortaliza'),
  @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44448, 4, 100112006, 'Repollo', 'Crespo record', 'Primera', 2000, 700, 800, 750, '$/unidad', 'Provincia del Elquí', 750, 1, 'Hortaliza'),
  @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44448, 4, 100112006, 'Repollo', 'Crespo record', 'Segunda', 1240, 550, 600, 575, '$/unidad', 'Provincia del Elquí', 575, 1, 'Hortaliza'),
  @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44448, 4, 100112006, 'Repollo', 'Morada(o)', 'Primera', 2000, 700, 800, 750, '$/unidad', 'Provincia del Elquí', 750, 1, 'Hortaliza'),
  @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', 44448, 4, 100112006, 'Repollo', 'Morada(o)', 'Segunda', 1240, 550, 600, 575, '$/unidad', 'Provincia del Elquí', 575, 1, 'Hortaliza')
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 1048 + $i
    $data = $newRows[$i]
    for ($c = 1; $c -le $data.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
